$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 22
    3  = 44
    4  = 67
    5  = 89
    6  = 111
    7  = 133
    8  = 156
    9  = 178
    10 = 200
    11 = 222
    12 = 244
    13 = 267
    14 = 289
    15 = 311
    16 = 333
    17 = 389
    18 = 444
    19 = 500
    20 = 556
    21 = 611
    22 = 667
    23 = 722
    24 = 778
    25 = 833
    26 = 889
    27 = 22
    28 = 44
    29 = 67
    30 = 89
    31 = 111
    32 = 133
    33 = 156
    34 = 178
    35 = 200
    36 = 222
    37 = 244
    38 = 267
    39 = 289
    40 = 311
    41 = 333
    42 = 389
    43 = 444
    44 = 500
    45 = 556
    46 = 611
    47 = 667
    48 = 722
    49 = 778
    50 = 833
    51 = 889
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $newValues[$row]
}

# Mirror column E's (best-fit) width onto column G, matching the author's
# formatting tweak alongside the data edits.
$ws.Range("G1").ColumnWidth = $ws.Range("E1").ColumnWidth

# Select the entire column F, matching the cursor position left in the
# saved workbook.
$ws.Range("F:F").Select() | Out-Null
